$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value2 = 34186351647.19485
$ws.Cells.Item(2, 5).Value2 = 39587899989.66781
$ws.Cells.Item(3, 4).Value2 = 37091099558.79767
$ws.Cells.Item(3, 5).Value2 = 42225084800.12243
$ws.Cells.Item(4, 4).Value2 = 40198076929.43095
$ws.Cells.Item(4, 5).Value2 = 45069517871.44103
$ws.Cells.Item(5, 4).Value2 = 43517672058.49415
$ws.Cells.Item(5, 5).Value2 = 48121821956.99799
$ws.Cells.Item(6, 4).Value2 = 47058396774.45271
$ws.Cells.Item(6, 5).Value2 = 51384483752.90667
$ws.Cells.Item(7, 4).Value2 = 50835693891.94943
$ws.Cells.Item(7, 5).Value2 = 54880987231.308
$ws.Cells.Item(8, 4).Value2 = 54872011449.47906
$ws.Cells.Item(8, 5).Value2 = 58655723627.69728
$ws.Cells.Item(9, 4).Value2 = 59182282313.27356
$ws.Cells.Item(9, 5).Value2 = 62714341389.51103
$ws.Cells.Item(10, 4).Value2 = 63786954273.91714
$ws.Cells.Item(10, 5).Value2 = 67071672214.66432
$ws.Cells.Item(11, 4).Value2 = 68713312825.13609
$ws.Cells.Item(11, 5).Value2 = 71760740127.67058
$ws.Cells.Item(12, 4).Value2 = 73984663655.60373
$ws.Cells.Item(12, 5).Value2 = 76804055532.77574
$ws.Cells.Item(13, 4).Value2 = 79628071072.34389
$ws.Cells.Item(13, 5).Value2 = 82227482280.71953
$ws.Cells.Item(14, 4).Value2 = 85673319223.69589
$ws.Cells.Item(14, 5).Value2 = 88060050493.15488
$ws.Cells.Item(15, 4).Value2 = 92152720738.72136
$ws.Cells.Item(15, 5).Value2 = 94339841643.58488
$ws.Cells.Item(16, 4).Value2 = 99102671017.52133
$ws.Cells.Item(16, 5).Value2 = 101100296037.0512
$ws.Cells.Item(17, 4).Value2 = 106557572271.5956
$ws.Cells.Item(17, 5).Value2 = 108373947245.5264
$ws.Cells.Item(18, 4).Value2 = 114557034788.3307
$ws.Cells.Item(18, 5).Value2 = 116193513804.62
$ws.Cells.Item(19, 4).Value2 = 123150275650.5164
$ws.Cells.Item(19, 5).Value2 = 124625705582.6505
$ws.Cells.Item(20, 4).Value2 = 132385025103.857
$ws.Cells.Item(20, 5).Value2 = 133709995342.9146
$ws.Cells.Item(21, 4).Value2 = 142306527916.3199
$ws.Cells.Item(21, 5).Value2 = 143491071241.3681
$ws.Cells.Item(22, 4).Value2 = 152965350222.193
$ws.Cells.Item(22, 5).Value2 = 154005487385.8299
$ws.Cells.Item(23, 4).Value2 = 164422745440.4424
$ws.Cells.Item(23, 5).Value2 = 165333635953.3872
$ws.Cells.Item(24, 4).Value2 = 176745571723.4986
$ws.Cells.Item(24, 5).Value2 = 177530903635.3451
$ws.Cells.Item(25, 4).Value2 = 190004377046.0385
$ws.Cells.Item(25, 5).Value2 = 190672402838.1364
$ws.Cells.Item(26, 4).Value2 = 204279776191.9457
$ws.Cells.Item(26, 5).Value2 = 204841598390.7104
$ws.Cells.Item(27, 4).Value2 = 219656802666.8462
$ws.Cells.Item(27, 5).Value2 = 220127408100.8791
$ws.Cells.Item(28, 4).Value2 = 236223744405.5191
$ws.Cells.Item(28, 5).Value2 = 236612379210.8285
$ws.Cells.Item(29, 4).Value2 = 254078847900.1587
$ws.Cells.Item(29, 5).Value2 = 254398593562.2359
$ws.Cells.Item(30, 4).Value2 = 273329107294.8923
$ws.Cells.Item(30, 5).Value2 = 273583882032.2896
$ws.Cells.Item(31, 4).Value2 = 294089135312.17
$ws.Cells.Item(31, 5).Value2 = 294294394642.3689
$ws.Cells.Item(32, 4).Value2 = 316479920861.035
$ws.Cells.Item(32, 5).Value2 = 316642220901.7303
$ws.Cells.Item(33, 4).Value2 = 340632354712.6805
$ws.Cells.Item(33, 5).Value2 = 340756723280.3546
$ws.Cells.Item(34, 4).Value2 = 366687459673.8379
$ws.Cells.Item(34, 5).Value2 = 366786645237.4137
$ws.Cells.Item(35, 4).Value2 = 394795500939.6024
$ws.Cells.Item(35, 5).Value2 = 394872984848.3019
$ws.Cells.Item(36, 4).Value2 = 425116985186.997
$ws.Cells.Item(36, 5).Value2 = 425176345700.1426
$ws.Cells.Item(37, 4).Value2 = 457825211061.2077
$ws.Cells.Item(37, 5).Value2 = 457869765162.8658
$ws.Cells.Item(38, 4).Value2 = 493107089999.0521
$ws.Cells.Item(38, 5).Value2 = 493139444093.6906
$ws.Cells.Item(39, 4).Value2 = 531163079873.5317
$ws.Cells.Item(39, 5).Value2 = 531186347602.708
$ws.Cells.Item(40, 4).Value2 = 572208938326.7451
$ws.Cells.Item(40, 5).Value2 = 572225277139.7828
$ws.Cells.Item(41, 4).Value2 = 616477438124.5095
$ws.Cells.Item(41, 5).Value2 = 616487839202.6716
$ws.Cells.Item(42, 4).Value2 = 664219538723.179
$ws.Cells.Item(42, 5).Value2 = 664226359652.0381
$ws.Cells.Item(43, 4).Value2 = 715705092344.9364
$ws.Cells.Item(43, 5).Value2 = 715709756445.9191
$ws.Cells.Item(44, 4).Value2 = 771224741288.8586
$ws.Cells.Item(44, 5).Value2 = 771227182473.2985
$ws.Cells.Item(45, 4).Value2 = 831092107181.0975
$ws.Cells.Item(45, 5).Value2 = 831093159499.7456
$ws.Cells.Item(46, 4).Value2 = 895645252769.3258
$ws.Cells.Item(46, 5).Value2 = 895645793122.1868
$ws.Cells.Item(47, 4).Value2 = 965248414276.1588
$ws.Cells.Item(47, 5).Value2 = 965248851162.6743
$ws.Cells.Item(48, 4).Value2 = 1040294097579.41
$ws.Cells.Item(48, 5).Value2 = 1040294482399.34
$ws.Cells.Item(49, 4).Value2 = 1121205412322.98
$ws.Cells.Item(49, 5).Value2 = 1121205751153.033
$ws.Cells.Item(50, 4).Value2 = 1208438481901.565
$ws.Cells.Item(50, 5).Value2 = 1208438777253.295
$ws.Cells.Item(51, 4).Value2 = 1302485012519.984
$ws.Cells.Item(51, 5).Value2 = 1302485266834.376
$ws.Cells.Item(52, 4).Value2 = 1403875058138.221
$ws.Cells.Item(52, 5).Value2 = 1403875274100.72
$ws.Cells.Item(53, 4).Value2 = 1513180006252.568
$ws.Cells.Item(53, 5).Value2 = 1513180186808.392
$ws.Cells.Item(54, 4).Value2 = 1631015800779.135
$ws.Cells.Item(54, 5).Value2 = 1631015949105.68
$ws.Cells.Item(55, 4).Value2 = 1758046413755.845
$ws.Cells.Item(55, 5).Value2 = 1758046533224.16
$ws.Cells.Item(56, 4).Value2 = 1894987584801.525
$ws.Cells.Item(56, 5).Value2 = 1894987678924.292
$ws.Cells.Item(57, 4).Value2 = 2042610853025.434
$ws.Cells.Item(57, 5).Value2 = 2042610925373.547
$ws.Cells.Item(58, 4).Value2 = 2201747908034.488
$ws.Cells.Item(58, 5).Value2 = 2201747962129.05
$ws.Cells.Item(59, 4).Value2 = 2373295282233.124
$ws.Cells.Item(59, 5).Value2 = 2373295321437.803
$ws.Cells.Item(60, 4).Value2 = 2558219408094.364
$ws.Cells.Item(60, 5).Value2 = 2558219435521.62
$ws.Cells.Item(61, 4).Value2 = 2757562068568.983
$ws.Cells.Item(61, 5).Value2 = 2757562087003.312
$ws.Cells.Item(62, 4).Value2 = 2972446272492.438
$ws.Cells.Item(62, 5).Value2 = 2972446284330.454
$ws.Cells.Item(63, 4).Value2 = 3204082590226.6
$ws.Cells.Item(63, 5).Value2 = 3204082597442.267
$ws.Cells.Item(64, 4).Value2 = 3453775983915.215
$ws.Cells.Item(64, 5).Value2 = 3453775988057.153
$ws.Cells.Item(65, 4).Value2 = 3722933168162.874
$ws.Cells.Item(65, 5).Value2 = 3722933170381.396
$ws.Cells.Item(66, 4).Value2 = 4013070542035.264
$ws.Cells.Item(66, 5).Value2 = 4013070543132.454
$ws.Cells.Item(67, 4).Value2 = 4325822738165.494
$ws.Cells.Item(67, 5).Value2 = 4325822738660.304
$ws.Cells.Item(68, 4).Value2 = 4662951837994.58
$ws.Cells.Item(68, 5).Value2 = 4662951838195.471
$ws.Cells.Item(69, 4).Value2 = 5026357304303.869
$ws.Cells.Item(69, 5).Value2 = 5026357304376.037
$ws.Cells.Item(70, 4).Value2 = 5418086686203.167
$ws.Cells.Item(70, 5).Value2 = 5418086686225.73
$ws.Cells.Item(71, 4).Value2 = 5840347157076.43
$ws.Cells.Item(71, 5).Value2 = 5840347157082.258
$ws.Cells.Item(72, 4).Value2 = 6295517951273.024
$ws.Cells.Item(72, 5).Value2 = 6295517951274.333
$ws.Cells.Item(73, 4).Value2 = 6786163770082.608
$ws.Cells.Item(73, 5).Value2 = 6786163770082.608
$ws.Cells.Item(74, 4).Value2 = 7315049232327.593
$ws.Cells.Item(74, 5).Value2 = 7315049232327.593
$ws.Cells.Item(75, 4).Value2 = 7885154450867.043
$ws.Cells.Item(75, 5).Value2 = 7885154450867.043
$ws.Cells.Item(76, 4).Value2 = 8499691823020.326
$ws.Cells.Item(76, 5).Value2 = 8499691823020.326
$ws.Cells.Item(77, 4).Value2 = 9162124129895.389
$ws.Cells.Item(77, 5).Value2 = 9162124129895.389
$ws.Cells.Item(78, 4).Value2 = 9876184046787.312
$ws.Cells.Item(78, 5).Value2 = 9876184046787.312
$ws.Cells.Item(79, 4).Value2 = 10645895174551.91
$ws.Cells.Item(79, 5).Value2 = 10645895174551.91
$ws.Cells.Item(80, 4).Value2 = 11475594710446.29
$ws.Cells.Item(80, 5).Value2 = 11475594710446.29
$ws.Cells.Item(81, 4).Value2 = 12369957886206.93
$ws.Cells.Item(81, 5).Value2 = 12369957886206.93
$ws.Cells.Item(82, 4).Value2 = 13334024311112.71
$ws.Cells.Item(82, 5).Value2 = 13334024311112.71
$ws.Cells.Item(83, 4).Value2 = 14373226368461.77
$ws.Cells.Item(83, 5).Value2 = 14373226368461.77
$ws.Cells.Item(84, 4).Value2 = 15493419825429.2
$ws.Cells.Item(84, 5).Value2 = 15493419825429.2
$ws.Cells.Item(85, 4).Value2 = 16700916828769.19
$ws.Cells.Item(85, 5).Value2 = 16700916828769.19
$ws.Cells.Item(86, 4).Value2 = 18002521472288.01
$ws.Cells.Item(86, 5).Value2 = 18002521472288.01
$ws.Cells.Item(87, 4).Value2 = 19405568136505.19
$ws.Cells.Item(87, 5).Value2 = 19405568136505.19
$ws.Cells.Item(88, 4).Value2 = 20917962816531.66
$ws.Cells.Item(88, 5).Value2 = 20917962816531.66
$ws.Cells.Item(89, 4).Value2 = 22548227671025.92
$ws.Cells.Item(89, 5).Value2 = 22548227671025.92
$ws.Cells.Item(90, 4).Value2 = 24305549043243.61
$ws.Cells.Item(90, 5).Value2 = 24305549043243.61
$ws.Cells.Item(91, 4).Value2 = 26199829224762.98
$ws.Cells.Item(91, 5).Value2 = 26199829224762.98
$ws.Cells.Item(92, 4).Value2 = 28241742253556.44
$ws.Cells.Item(92, 5).Value2 = 28241742253556.44
$ws.Cells.Item(93, 4).Value2 = 30442794060808.52
$ws.Cells.Item(93, 5).Value2 = 30442794060808.52
$ws.Cells.Item(94, 4).Value2 = 32815387305386.0
$ws.Cells.Item(94, 5).Value2 = 32815387305386.0
$ws.Cells.Item(95, 4).Value2 = 35372891261284.34
$ws.Cells.Item(95, 5).Value2 = 35372891261284.34
$ws.Cells.Item(96, 4).Value2 = 38129717151847.3
$ws.Cells.Item(96, 5).Value2 = 38129717151847.3
$ws.Cells.Item(97, 4).Value2 = 41101399355250.23
$ws.Cells.Item(97, 5).Value2 = 41101399355250.23
